# Scheduled-runner refresh of Universalis market-board pulls.
# Updates currentAveragePrice(NQ/HQ) and the derived Leve price / profit
# columns (H,I,J,K,L,M,N) on each job sheet to the latest snapshot.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 298.25
$ws.Range("I2").Value = 298.25
$ws.Range("K2").Value = 298.25
$ws.Range("M2").Value = -185.25
$ws.Range("H32").Value = 1125
$ws.Range("J32").Value = 1125
$ws.Range("L32").Value = 1125
$ws.Range("N32").Value = -1777
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()
$ws.Range("H54").Value = 9800
$ws.Range("J54").Value = 9800
$ws.Range("L54").Value = 9800
$ws.Range("N54").Value = -10772
$ws.Range("H70").Value = 1875
$ws.Range("J70").Value = 1714.2858
$ws.Range("L70").Value = 5142.857400000001
$ws.Range("N70").Value = -5682.857400000001
$ws.Range("H73").Value = 1875
$ws.Range("J73").Value = 1714.2858
$ws.Range("L73").Value = 5142.857400000001
$ws.Range("N73").Value = -7014.857400000001
$ws.Range("H86").Value = 9000
$ws.Range("I86").Value = 9000
$ws.Range("K86").Value = 9000
$ws.Range("M86").Value = -7877
$ws.Range("H89").Value = 9000
$ws.Range("I89").Value = 9000
$ws.Range("K89").Value = 45000
$ws.Range("M89").Value = -39384
$ws.Range("H132").Value = 1953.3334
$ws.Range("I132").Value = 947.5625
$ws.Range("K132").Value = 2842.6875
$ws.Range("M132").Value = -312.6875
$ws.Range("H135").Value = 707.875
$ws.Range("I135").Value = 666.1429000000001
$ws.Range("K135").Value = 5995.2861
$ws.Range("M135").Value = -3460.2861
$ws.Range("H138").Value = 12076.762
$ws.Range("J138").Value = 12821.711
$ws.Range("L138").Value = 38465.133
$ws.Range("N138").Value = -48745.133

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1211.05
$ws.Range("I2").Value = 1171.375
$ws.Range("K2").Value = 1171.375
$ws.Range("M2").Value = -1058.375
$ws.Range("H62").Value = 10000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 10000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 10000
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -11248
$ws.Range("H65").Value = 10000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 10000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 30000
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -36240
$ws.Range("H110").Value = 1353.6666
$ws.Range("I110").Value = 1345.0588
$ws.Range("J110").Value = 1500
$ws.Range("K110").Value = 1345.0588
$ws.Range("L110").Value = 1500
$ws.Range("M110").Value = 699.9412
$ws.Range("N110").Value = -5590
$ws.Range("H116").Value = 1211.05
$ws.Range("I116").Value = 1171.375
$ws.Range("K116").Value = 1171.375
$ws.Range("M116").Value = 1122.625
$ws.Range("H122").Value = 2143.8572
$ws.Range("I122").Value = 2084.5
$ws.Range("K122").Value = 6253.5
$ws.Range("M122").Value = -3803.5

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1211.05
$ws.Range("I3").Value = 1171.375
$ws.Range("K3").Value = 1171.375
$ws.Range("M3").Value = -1057.375
$ws.Range("H20").Value = 2500
$ws.Range("I20").Value = 2500
$ws.Range("K20").Value = 2500
$ws.Range("M20").Value = -2253
$ws.Range("H57").Value = 80000
$ws.Range("J57").Value = 80000
$ws.Range("L57").Value = 80000
$ws.Range("N57").Value = -81440
$ws.Range("H86").Value = 6221.9443
$ws.Range("I86").Value = 2521.889
$ws.Range("K86").Value = 2521.889
$ws.Range("M86").Value = -1398.889
$ws.Range("H89").Value = 6221.9443
$ws.Range("I89").Value = 2521.889
$ws.Range("K89").Value = 12609.445
$ws.Range("M89").Value = -6993.445
$ws.Range("H102").Value = 30556
$ws.Range("I102").Value = 30556
$ws.Range("K102").Value = 30556
$ws.Range("M102").Value = -27311
$ws.Range("H136").Value = 80000
$ws.Range("J136").Value = 80000
$ws.Range("L136").Value = 80000
$ws.Range("N136").Value = -90200

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 48039.5
$ws.Range("I59").Value = 4904
$ws.Range("J59").Value = 62418
$ws.Range("K59").Value = 4904
$ws.Range("L59").Value = 62418
$ws.Range("M59").Value = -3759
$ws.Range("N59").Value = -64708
$ws.Range("H60").Value = 35818.6
$ws.Range("J60").Value = 49699.332
$ws.Range("L60").Value = 49699.332
$ws.Range("N60").Value = -50721.332

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 1639.8
$ws.Range("J29").Value = 1639.8
$ws.Range("L29").Value = 4919.4
$ws.Range("N29").Value = -5473.4
$ws.Range("H131").Value = 943
$ws.Range("I131").Value = 863.4286
$ws.Range("K131").Value = 2590.2858
$ws.Range("M131").Value = 2449.7142

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 12124.6875
$ws.Range("I22").Value = 15455
$ws.Range("J22").Value = 7842.857
$ws.Range("K22").Value = 15455
$ws.Range("L22").Value = 7842.857
$ws.Range("M22").Value = -15160
$ws.Range("N22").Value = -8432.857
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("M24").ClearContents()
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()
$ws.Range("H27").Value = 12124.6875
$ws.Range("I27").Value = 15455
$ws.Range("J27").Value = 7842.857
$ws.Range("K27").Value = 15455
$ws.Range("L27").Value = 7842.857
$ws.Range("M27").Value = -15348
$ws.Range("N27").Value = -8056.857
$ws.Range("H41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()
$ws.Range("H55").Value = 692.6667
$ws.Range("I55").Value = 692.6667
$ws.Range("K55").Value = 692.6667
$ws.Range("M55").Value = -519.6667
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()
$ws.Range("H136").Value = 2414.2856
$ws.Range("I136").Value = 2414.2856
$ws.Range("K136").Value = 7242.8568
$ws.Range("M136").Value = -4692.8568
